# Updates Siglec1-Spn TPM results: rows 2-11 get new recomputed values and
# rows 12-16 are newly appended (dimension grows from A1:T11 to A1:T16).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Siglec1"
$ws.Range("C2").Value = "Spn"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.093736
$ws.Range("H2").Value = 3.281208
$ws.Range("I2").Value = 0.004351073087712599
$ws.Range("J2").Value = 0.004369500290106573
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.359683333333334
$ws.Range("N2").Value = 16.07905
$ws.Range("O2").Value = 0.5495559766256753
$ws.Range("P2").Value = 0.5567946816040513
$ws.Range("Q2").Value = 5.862078610266668
$ws.Range("R2").Value = 52.75870749240001
$ws.Range("S2").Value = 0.00239115822008759
$ws.Range("T2").Value = 0.002432914522798699

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Siglec1"
$ws.Range("C3").Value = "Spn"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.093736
$ws.Range("H3").Value = 3.281208
$ws.Range("I3").Value = 0.004351073087712599
$ws.Range("J3").Value = 0.004369500290106573
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.380377
$ws.Range("N3").Value = 0.760754
$ws.Range("O3").Value = 0.03900201573866823
$ws.Range("P3").Value = 0.02634383133387908
$ws.Range("Q3").Value = 0.4160320184720001
$ws.Range("R3").Value = 2.496192110832
$ws.Range("S3").Value = 0.0001697006210470626
$ws.Range("T3").Value = 0.0001151093786559033

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Siglec1"
$ws.Range("C4").Value = "Spn"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.093736
$ws.Range("H4").Value = 3.281208
$ws.Range("I4").Value = 0.004351073087712599
$ws.Range("J4").Value = 0.004369500290106573
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.012692
$ws.Range("N4").Value = 12.038076
$ws.Range("O4").Value = 0.4114420076356565
$ws.Range("P4").Value = 0.4168614870620697
$ws.Range("Q4").Value = 4.388825697312001
$ws.Range("R4").Value = 39.499431275808
$ws.Range("S4").Value = 0.001790214246577947
$ws.Range("T4").Value = 0.001821476388651971

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Siglec1"
$ws.Range("C5").Value = "Spn"
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4782236666666666
$ws.Range("H5").Value = 1.434671
$ws.Range("I5").Value = 0.001902457380885857
$ws.Range("J5").Value = 0.001910514466229354
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.359683333333334
$ws.Range("N5").Value = 16.07905
$ws.Range("O5").Value = 0.5495559766256753
$ws.Range("P5").Value = 0.5567946816040513
$ws.Range("Q5").Value = 2.563127415838889
$ws.Range("R5").Value = 23.06814674255
$ws.Range("S5").Value = 0.001045506823941452
$ws.Range("T5").Value = 0.001063764293924107

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Siglec1"
$ws.Range("C6").Value = "Spn"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4782236666666666
$ws.Range("H6").Value = 1.434671
$ws.Range("I6").Value = 0.001902457380885857
$ws.Range("J6").Value = 0.001910514466229354
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.380377
$ws.Range("N6").Value = 0.760754
$ws.Range("O6").Value = 0.03900201573866823
$ws.Range("P6").Value = 0.02634383133387908
$ws.Range("Q6").Value = 0.1819052836556667
$ws.Range("R6").Value = 1.091431701934
$ws.Range("S6").Value = 0.00007419967271145574
$ws.Range("T6").Value = 0.00005033027085928213

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Siglec1"
$ws.Range("C7").Value = "Spn"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4782236666666666
$ws.Range("H7").Value = 1.434671
$ws.Range("I7").Value = 0.001902457380885857
$ws.Range("J7").Value = 0.001910514466229354
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.012692
$ws.Range("N7").Value = 12.038076
$ws.Range("O7").Value = 0.4114420076356565
$ws.Range("P7").Value = 0.4168614870620697
$ws.Range("Q7").Value = 1.918964281444
$ws.Range("R7").Value = 17.270678532996
$ws.Range("S7").Value = 0.0007827508842329498
$ws.Range("T7").Value = 0.0007964199014459649

# Row 8
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Siglec1"
$ws.Range("C8").Value = "Spn"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 125.3706816666667
$ws.Range("H8").Value = 376.112045
$ws.Range("I8").Value = 0.4987464973156379
$ws.Range("J8").Value = 0.5008587354840279
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.359683333333334
$ws.Range("N8").Value = 16.07905
$ws.Range("O8").Value = 0.5495559766256753
$ws.Range("P8").Value = 0.5567946816040513
$ws.Range("Q8").Value = 671.9471530174723
$ws.Range("R8").Value = 6047.52437715725
$ws.Range("S8").Value = 0.2740891184209301
$ws.Range("T8").Value = 0.2788754801524371

# Row 9
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Siglec1"
$ws.Range("C9").Value = "Spn"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 125.3706816666667
$ws.Range("H9").Value = 376.112045
$ws.Range("I9").Value = 0.4987464973156379
$ws.Range("J9").Value = 0.5008587354840279
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.5
$ws.Range("M9").Value = 0.380377
$ws.Range("N9").Value = 0.760754
$ws.Range("O9").Value = 0.03900201573866823
$ws.Range("P9").Value = 0.02634383133387908
$ws.Range("Q9").Value = 47.68812378032167
$ws.Range("R9").Value = 286.12874268193
$ws.Range("S9").Value = 0.01945211873791016
$ws.Range("T9").Value = 0.01319453804969119

# Row 10
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Siglec1"
$ws.Range("C10").Value = "Spn"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 125.3706816666667
$ws.Range("H10").Value = 376.112045
$ws.Range("I10").Value = 0.4987464973156379
$ws.Range("J10").Value = 0.5008587354840279
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.012692
$ws.Range("N10").Value = 12.038076
$ws.Range("O10").Value = 0.4114420076356565
$ws.Range("P10").Value = 0.4168614870620697
$ws.Range("Q10").Value = 503.07393135838
$ws.Range("R10").Value = 4527.665382225419
$ws.Range("S10").Value = 0.2052052601567976
$ws.Range("T10").Value = 0.2087887172818997

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Siglec1"
$ws.Range("C11").Value = "Spn"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.1802775
$ws.Range("H11").Value = 6.360555
$ws.Range("I11").Value = 0.01265170008275114
$ws.Range("J11").Value = 0.008470187479043942
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.359683333333334
$ws.Range("N11").Value = 16.07905
$ws.Range("O11").Value = 0.5495559766256753
$ws.Range("P11").Value = 0.5567946816040513
$ws.Range("Q11").Value = 17.045280312125
$ws.Range("R11").Value = 102.27168187275
$ws.Range("S11").Value = 0.006952817394951442
$ws.Range("T11").Value = 0.004716155340520894

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Siglec1"
$ws.Range("C12").Value = "Spn"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.1802775
$ws.Range("H12").Value = 6.360555
$ws.Range("I12").Value = 0.01265170008275114
$ws.Range("J12").Value = 0.008470187479043942
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 0.380377
$ws.Range("N12").Value = 0.760754
$ws.Range("O12").Value = 0.03900201573866823
$ws.Range("P12").Value = 0.02634383133387908
$ws.Range("Q12").Value = 1.2097044146175
$ws.Range("R12").Value = 4.83881765847
$ws.Range("S12").Value = 0.0004934418057483702
$ws.Range("T12").Value = 0.000223137190314268

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Siglec1"
$ws.Range("C13").Value = "Spn"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.1802775
$ws.Range("H13").Value = 6.360555
$ws.Range("I13").Value = 0.01265170008275114
$ws.Range("J13").Value = 0.008470187479043942
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 4.012692
$ws.Range("N13").Value = 12.038076
$ws.Range("O13").Value = 0.4114420076356565
$ws.Range("P13").Value = 0.4168614870620697
$ws.Range("Q13").Value = 12.76147408203
$ws.Range("R13").Value = 76.56884449217999
$ws.Range("S13").Value = 0.005205440882051332
$ws.Range("T13").Value = 0.003530894948208781

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Siglec1"
$ws.Range("C14").Value = "Spn"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 121.2486343333333
$ws.Range("H14").Value = 363.745903
$ws.Range("I14").Value = 0.4823482721330124
$ws.Range("J14").Value = 0.4843910622805923
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 5.359683333333334
$ws.Range("N14").Value = 16.07905
$ws.Range("O14").Value = 0.5495559766256753
$ws.Range("P14").Value = 0.5567946816040513
$ws.Range("Q14").Value = 649.8542846257944
$ws.Range("R14").Value = 5848.688561632151
$ws.Range("S14").Value = 0.2650773757657647
$ws.Range("T14").Value = 0.2697063672943706

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Siglec1"
$ws.Range("C15").Value = "Spn"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 121.2486343333333
$ws.Range("H15").Value = 363.745903
$ws.Range("I15").Value = 0.4823482721330124
$ws.Range("J15").Value = 0.4843910622805923
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.5
$ws.Range("M15").Value = 0.380377
$ws.Range("N15").Value = 0.760754
$ws.Range("O15").Value = 0.03900201573866823
$ws.Range("P15").Value = 0.02634383133387908
$ws.Range("Q15").Value = 46.12019178181033
$ws.Range("R15").Value = 276.721150690862
$ws.Range("S15").Value = 0.01881255490125118
$ws.Range("T15").Value = 0.01276071644435844

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Siglec1"
$ws.Range("C16").Value = "Spn"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 121.2486343333333
$ws.Range("H16").Value = 363.745903
$ws.Range("I16").Value = 0.4823482721330124
$ws.Range("J16").Value = 0.4843910622805923
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 4.012692
$ws.Range("N16").Value = 12.038076
$ws.Range("O16").Value = 0.4114420076356565
$ws.Range("P16").Value = 0.4168614870620697
$ws.Range("Q16").Value = 486.533425000292
$ws.Range("R16").Value = 4378.800825002628
$ws.Range("S16").Value = 0.1984583414659966
$ws.Range("T16").Value = 0.2019239785418633

